$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-style row 11 (it moves from being the last data row to a "middle"
#    row of the table, so it switches from the s4/s5 style pair to s8/s9,
#    matching the look of rows 4-8).
# ---------------------------------------------------------------------------
$ws.Range("A4:B4").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("C4:E4").Copy()
$ws.Range("C11:E11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Row 12 - new data row (SCRIPT/P01P04A/us0402.ssb)
#    Shared-string insertion order must be C, A, D, E to line up with the
#    workbook being produced.
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = " This is the first time I\'ve been\nhere, but it\'s a pretty fun place."
$ws.Range("A12").Value = "SCRIPT/P01P04A/us0402.ssb"
$ws.Range("D12").Value = " Я здесь впервые, но мне уже\nнравится это место."
$ws.Range("E12").Value = " Ÿ èäåòû âðåñâúå, îï íîå ôçå\nîñàâéóòÿ üóï íåòóï."
$ws.Range("B12").Value = 57

$ws.Range("A2:B2").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C12:E12").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = 48.6

# ---------------------------------------------------------------------------
# 3. Row 13 - group header row (only column A has text, rest stay blank,
#    matching the look of row 10).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "SCRIPT/P01P04A/us2007.ssb"

$ws.Range("A10:E10").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Rows.Item(13).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 4. Row 14 - new data row (SCRIPT/T01P02A/us2009.ssb )
#    Shared-string insertion order is C, A, D, E again.
# ---------------------------------------------------------------------------
$ws.Range("C14").Value = " The rumor is that [CS:N]Drowzee[CR] has\nreformed and forsaken his criminal past.[K]\nBut I have to wonder…"
$ws.Range("A14").Value = "SCRIPT/T01P02A/us2009.ssb "
$ws.Range("D14").Value = " Ходят слухи, что [CS:N]Дроузи[CR]\nисправился и искупил своё криминальное\nпрошлое.[K] Но я мало этому верю..."
$ws.Range("E14").Value = " Öïäÿó òìôöé, œóï [CS:N]Äñïôèé[CR]\néòðñàâéìòÿ é éòëôðéì òâïæ ëñéíéîàìûîïå\nðñïšìïå.[K] Îï ÿ íàìï üóïíô âåñý..."
$ws.Range("B14").Value = 38

$ws.Range("A4:B4").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("C4:E4").Copy()
$ws.Range("C14:E14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 5. Row 15 - new data row (SCRIPT/T01P02A/us2013.ssb)
#    Shared-string insertion order is A, C, D, E this time.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "SCRIPT/T01P02A/us2013.ssb"
$ws.Range("C15").Value = " I hear that [CS:N]Drowzee[CR] is at\nthe guild."
$ws.Range("D15").Value = " Говорят, что [CS:N]Дроузи[CR] находится\nв гильдии."
$ws.Range("E15").Value = " Ãïâïñÿó, œóï [CS:N]Äñïôèé[CR] îàöïäéóòÿ\nâ ãéìûäéé."
$ws.Range("B15").Value = 18

$ws.Range("A2:B2").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("C2:E2").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 6. Update the visible selection to match the new bottom of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("D15").Select()
